# "changed internal battery measurement"
# The workbook's Tabelle1 sheet computes a battery-factor calculation.
# F10 holds the measured "Return from ADC" value that the rest of the
# sheet (F13, F16, H17) derives via formulas. Updating it lets Excel's
# normal recalculation engine refresh all of the dependent cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F10").Value = 266

# Leave the selection where the author left it when they saved the file.
$ws.Range("F11").Select()
